$d = $word.ActiveDocument

$d.Content.Find.Execute("disposition_method", $true, $false, $false, $false, $false,
                         $true, 1, $false, "disposition", 2)
